# edit.ps1 - applies the documented diff to documentation.docx
#
# Summary of changes:
#  1) "Tiedot noudettu lukuvuoden 2020-2021 tiedoista ..." paragraph: the
#     year "2020-2021" was wrapped in proofErr gramStart/gramEnd markers
#     across three runs; collapse back into a single plain run.
#  2) The (first) "Sovittu:" / "Toteutunut:" pair that follows "Anni" gets
#     "Sovittu:" turned bold, and two new bullet-like paragraphs plus a
#     blank paragraph are inserted between "Sovittu:" and "Toteutunut:".
#  3) "Tutkinto-ohjelmien rakenteen lukeminen Kori APIsta" paragraph gains
#     a <w:lastRenderedPageBreak/> marker at the start of its first run.
#  4) "Luokkakaavion tekeminen dokumentointia varten" paragraph loses the
#     <w:lastRenderedPageBreak/> marker it used to carry (it moved to #3).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge the "2020-2021" runs (drops the gramStart/gramEnd proofErr
#    pair) back into a single run with the full sentence.
# ---------------------------------------------------------------------
$find = $d.Content.Find
$find.Execute(
    "Tiedot noudettu lukuvuoden 2020-2021 tiedoista (ohjeistuksessa annettu linkki).",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Tiedot noudettu lukuvuoden 2020-2021 tiedoista (ohjeistuksessa annettu linkki).",
    2
) | Out-Null

# ---------------------------------------------------------------------
# 2) Locate the "Sovittu:" paragraph that immediately follows the "Anni"
#    heading paragraph (the first Sovittu:/Toteutunut: pair in the doc).
# ---------------------------------------------------------------------
$sovittuIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq "Anni") {
        $sovittuIndex = $i + 1
        break
    }
}

$sovittuRange = $d.Paragraphs.Item($sovittuIndex).Range
$sovittuRange.Collapse(0)
$sovittuRange.InsertParagraphAfter() | Out-Null
$sovittuRange.InsertParagraphAfter() | Out-Null
$sovittuRange.InsertParagraphAfter() | Out-Null

# New paragraph: "Testien luominen luokalle StudyTree"
$xmlWrap = @'
<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>{0}</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$studyTreeParaXml = $xmlWrap -f @'
<w:p><w:r><w:t xml:space="preserve">Testien luominen luokalle </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>StudyTree</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$d.Paragraphs.Item($sovittuIndex + 1).Range.InsertXML($studyTreeParaXml) | Out-Null

# New paragraph: "Javadoc-kommenttien luominen niille metodeille ja luokille joista ne puuttuvat"
$javadocParaXml = $xmlWrap -f @'
<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Javadoc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-kommenttien luominen niille metodeille ja luokille joista ne puuttuvat</w:t></w:r></w:p>
'@
$d.Paragraphs.Item($sovittuIndex + 2).Range.InsertXML($javadocParaXml) | Out-Null

# New blank paragraph between the inserted text and "Toteutunut:"
$blankParaXml = $xmlWrap -f '<w:p/>'
$d.Paragraphs.Item($sovittuIndex + 3).Range.InsertXML($blankParaXml) | Out-Null

# Make the original "Sovittu:" paragraph bold (both pPr/rPr and run rPr),
# done last so the freshly-inserted blank paragraphs above don't inherit
# the bold paragraph-mark formatting.
$boldSovittuXml = $xmlWrap -f @'
<w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Sovittu:</w:t></w:r></w:p>
'@
$d.Paragraphs.Item($sovittuIndex).Range.InsertXML($boldSovittuXml) | Out-Null

# ---------------------------------------------------------------------
# 3) & 4) Move <w:lastRenderedPageBreak/> from the "Luokkakaavion..."
#    paragraph to the "Tutkinto-ohjelmien rakenteen lukeminen..." one.
# ---------------------------------------------------------------------
$tutkintoIndex = 0
$luokkaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)
    if ($t -eq "Tutkinto-ohjelmien rakenteen lukeminen Kori APIsta") {
        $tutkintoIndex = $i
    }
    if ($t -eq "Luokkakaavion tekeminen dokumentointia varten") {
        $luokkaIndex = $i
    }
}

$tutkintoXml = $xmlWrap -f @'
<w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Tutkinto-ohjelmien rakenteen lukeminen Kori </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>APIsta</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$d.Paragraphs.Item($tutkintoIndex).Range.InsertXML($tutkintoXml) | Out-Null

$luokkaXml = $xmlWrap -f @'
<w:p><w:r><w:t>Luokkakaavion tekeminen dokumentointia varten</w:t></w:r></w:p>
'@
$d.Paragraphs.Item($luokkaIndex).Range.InsertXML($luokkaXml) | Out-Null
